$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from Ferlab.bio CodeS" sheet to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$includeSheet.Name = "Include #0"

# 2. Update the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")

# Update Date value (row 8, column B)
$meta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Update Contact value (row 10, column B)
$meta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row after row 10 (Contact) for "Jurisdiction"
$meta.Rows.Item(11).Insert()

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("A11").Style = $meta.Range("A10").Style
$meta.Range("B11").Value = ""
$meta.Range("B11").Style = $meta.Range("B10").Style
